# Allow alt_cells for the relative formula type.
#
# This reproduces the authoring change made to MasterGalleryKey.xlsx:
#  - On the "Relative Samples" sheet, the old "Complex formula" placeholder
#    row (row 5, only cell A5 filled in) is replaced by two new worked
#    examples that mirror the existing A2/A3 rows (A = B + C):
#      * row 4 demonstrates the "delta" rubric tolerance
#      * row 5 demonstrates the new "alt_cells" rubric tolerance, which
#        points at F10 as an alternate acceptable cell
#  - Matching grading-rubric comments are added on A4 and A5.
#  - A new helper value (F10 = 8888) is added as the alt_cells target.
#  - The "Relative Samples_CheckOrder" sheet is updated so the grading
#    order list includes the new A4/A5 cells (and keeps A6 after them).
#  - A few sheet-view selections shift around as a side effect of the
#    authors having clicked around while editing.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Relative Samples" sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Relative Samples")

# New row 4: same shape as rows 2/3 (A = B + C), with a note that this
# cell is meant to exercise the "delta" tolerance.
$ws.Range("A4").Formula = "=B4+C4"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = "// delta for being more tolerant"

# New row 5 replaces the old "Complex formula" placeholder (A5 only).
# Same shape again, but the note now calls out the "alt_cells" tolerance.
$ws.Range("A5").ClearContents()
$ws.Range("A5").Formula = "=B5+C5"
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 200
$ws.Range("D5").Value = "// alt cells for being more tolerant"

# Helper cell used as the alternate acceptable answer for A5's rubric.
$ws.Range("F10").Value = 8888

# Grading-rubric comments for the two new cells.
$a4Comment = "rubric:`n score: 1`n type: relative`n delta: 5"
$ws.Range("A4").AddComment($a4Comment) | Out-Null

$a5Comment = "rubric:`n score: 1`n type: relative`nalt_cells:`n - F10"
$ws.Range("A5").AddComment($a5Comment) | Out-Null

# ----------------------------------------------------------------------
# "Relative Samples_CheckOrder" sheet
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Relative Samples_CheckOrder")

# Row 4 used to point at "A6" (the only graded cell after A2/A3); now
# that A4 has a rubric too, row 4 points at "A4" instead.
$ws2.Range("B4").Value = "A4"

# New row 5 points at the new "A5" rubric cell.
$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "A5"

# New row 6 keeps "A6" in the grading order, now listed after A4/A5.
$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = "A6"

# ----------------------------------------------------------------------
# Cosmetic sheet-view selection changes
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Constant Samples")
$ws4.Range("B9").Select()

$ws2.Select()
$ws2.Range("D10").Select()

$ws.Activate()
$ws.Range("E12").Select()

Write-Host "Applied: alt_cells support for relative formula type"
